# Insert two new data rows into the Acelga (Chillán) price sheet.
# The new rows are inserted immediately above the current row 467, pushing
# the existing rows 467:566 down to 469:568 (dimension grows from R566 to R568).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before row 467 (formatting/row-height is copied from
# the row above by default, which also gives column D the date style it
# needs).
$ws.Rows.Item(467).Resize(2).Insert()

# ---- New row 467 ------------------------------------------------------
$ws.Cells.Item(467, 1).Value  = 7
$ws.Cells.Item(467, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(467, 3).Value  = "Ñuble"
$ws.Cells.Item(467, 4).Value  = 45173
$ws.Cells.Item(467, 5).Value  = 16
$ws.Cells.Item(467, 6).Value  = 100112009
$ws.Cells.Item(467, 7).Value  = "Acelga"
$ws.Cells.Item(467, 8).Value  = "Sin especificar"
$ws.Cells.Item(467, 9).Value  = "Primera"
$ws.Cells.Item(467, 10).Value = 250
$ws.Cells.Item(467, 11).Value = 700
$ws.Cells.Item(467, 12).Value = 700
$ws.Cells.Item(467, 13).Value = 700
$ws.Cells.Item(467, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(467, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(467, 16).Value = 700
$ws.Cells.Item(467, 17).Value = 1
$ws.Cells.Item(467, 18).Value = "Hortaliza"

# ---- New row 468 ------------------------------------------------------
$ws.Cells.Item(468, 1).Value  = 7
$ws.Cells.Item(468, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(468, 3).Value  = "Ñuble"
$ws.Cells.Item(468, 4).Value  = 45173
$ws.Cells.Item(468, 5).Value  = 16
$ws.Cells.Item(468, 6).Value  = 100112009
$ws.Cells.Item(468, 7).Value  = "Acelga"
$ws.Cells.Item(468, 8).Value  = "Sin especificar"
$ws.Cells.Item(468, 9).Value  = "Segunda"
$ws.Cells.Item(468, 10).Value = 200
$ws.Cells.Item(468, 11).Value = 500
$ws.Cells.Item(468, 12).Value = 500
$ws.Cells.Item(468, 13).Value = 500
$ws.Cells.Item(468, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(468, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(468, 16).Value = 500
$ws.Cells.Item(468, 17).Value = 1
$ws.Cells.Item(468, 18).Value = "Hortaliza"
